$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Hora (G) columns are stored as text in the workbook.
# Apply a Text number format first so the new values keep their text type
# (matching the original inlineStr cells) instead of becoming numbers.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Update Price (column D) values that changed
$ws.Range("D2").Value = "269.61"
$ws.Range("D3").Value = "22.90"
$ws.Range("D4").Value = "6.372"
$ws.Range("D5").Value = "0.06240"
$ws.Range("D6").Value = "3.649"
$ws.Range("D7").Value = "6.700"
$ws.Range("D8").Value = "1.372"
$ws.Range("D9").Value = "0.8334"
$ws.Range("D10").Value = "0.01377"
$ws.Range("D11").Value = "0.1632"
$ws.Range("D12").Value = "0.08393"
$ws.Range("D13").Value = "0.03464"
$ws.Range("D14").Value = "0.03130"
$ws.Range("D15").Value = "0.09321"
$ws.Range("D16").Value = "3.882"
$ws.Range("D17").Value = "0.001706"
$ws.Range("D18").Value = "0.04813"
$ws.Range("D19").Value = "0.006242"
$ws.Range("D20").Value = "0.001087"
$ws.Range("D21").Value = "0.003452"
$ws.Range("D23").Value = "3.736"
$ws.Range("D24").Value = "2.369"
$ws.Range("D25").Value = "0.3405"
$ws.Range("D40").Value = "0.04684"
$ws.Range("D41").Value = "0.006900"
$ws.Range("D42").Value = "0.1167"
$ws.Range("D43").Value = "0.003311"
$ws.Range("D44").Value = "0.01119"
$ws.Range("D45").Value = "0.00006271"
$ws.Range("D47").Value = "0.8798"
$ws.Range("D48").Value = "0.08049"

# Update Hora (column G) from 13 to 14 for every data row
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 7).Value = "14"
}
